# Generate Report for Handback
# Update the "generated/xliff" timestamps recorded on each sheet to reflect
# a fresh report generation run.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# in-sync-with-en-US row (mirrors the de-de Correspond Handoff Datetime
# for the same source file).
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 03:04:34"

# zh-cn sheet: Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K) for the first data row.
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 03:04:27"
$wsZhCn.Range("K2").Value = "2016-08-24 03:04:58"

# de-de sheet: Correspond Handoff Datetime (H) and
# Correspond Handback DateTime (K) for the first data row.
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 03:04:34"
$wsDeDe.Range("K2").Value = "2016-08-24 03:05:12"
